$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vimrcText = @'
" myles put 2017-05-09 on indent tab to space
filetype plugin indent on       
set tabstop=4                   "show existing tab with 4 spaces width
set shiftwidth=4                " when indenting with '>', use 4 spaces width
set expandtab                   " On pressing tab, insert 4 spaces
" myles put 2017-05-09 on set default colorscheme
colorscheme darkblue
" Install pathogen 2017-05-23
execute pathogen#infect()
syntax on 
" 2017-05-23 set option
"set nowrap
" 2017-05-24 set option
set incsearch
" 2017-05-24 Open Nerdtree when vim open with no input
"autocmd StdinReadPre * let s:std_in=1
"autocmd VimEnter * if argc() == 0 && !exists("s:std_in") | NERDTree | endif
" 2017-05-24 Create shortcut for NERDTree
map <C-n> :NERDTreeToggle<CR>
" 2017-05-25 Change the directory arrow
let g:NERDTreeDirArrowExpandable="+"
let g:NERDTreeDirArrowCollapsible="-"
" 2017-05-25 Change the directory arrow
autocmd BufWinLeave *.* mkview
"autocmd BufWinEnter *.* silent loadview
" 2017-05-25 Syntastic plugin 
autocmd BufWinEnter *.* silent SyntasticToggleMode
set statusline+=%#warningmsg#
set statusline+=%{SyntasticStatuslineFlag()}
set statusline+=%*
let g:syntastic_always_populate_loc_list = 1
let g:syntastic_auto_loc_list = 1
"let g:syntastic_check_on_open = 1
"let g:syntastic_check_on_wq = 0
" 2017-05-26 Syntax folding autostart
set foldmethod=syntax
set foldlevel=1
'@

$ws.Range("A28").Value = "vimrc"
$ws.Range("B28").Value = "snapshot"
$ws.Range("C28").Value = $vimrcText

$ws.Rows.Item(28).RowHeight = 92.25

$ws.Application.ActiveWindow.ScrollRow = 27
$ws.Range("C34").Select()
